$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = "Ali Al Memari"
$summary.Range("B4").Value = 1579.24
$summary.Range("B6").Value = 2875
$summary.Range("B7").Value = 37121
$summary.Range("B8").Value = -34246
$summary.Range("B9").Value = 0.08

# ---------------------------------------------------------------------
# Assets sheet - remove both "Vehicles" rows (Luxury Car, Premium Car),
# keep the Liquid Assets / Savings Account row and the TOTAL ASSETS row,
# then update the remaining numeric values.
# ---------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")
$assets.Rows.Item(2).Delete()
$assets.Rows.Item(2).Delete()
$assets.Range("C2").Value = 2875
$assets.Range("C3").Value = 2875

# ---------------------------------------------------------------------
# Liabilities sheet - remove the "Auto Loans" row, keep the Credit Cards
# row and the TOTAL LIABILITIES row, then update the remaining values.
# ---------------------------------------------------------------------
$liabilities = $wb.Worksheets.Item("Liabilities")
$liabilities.Rows.Item(2).Delete()
$liabilities.Range("C2").Value = 37121
$liabilities.Range("D2").Value = 1856
$liabilities.Range("C3").Value = 37121
